$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row's double rule (its own bottom border sitting right on top
# of row 4's top border) gets simplified down to a single line now that
# the table is being extended - drop the redundant bottom border under the
# year headers.
$ws.Range("B3:J3").Borders.Item(9).LineStyle = -4142  # xlEdgeBottom, xlLineStyleNone

# Copy formatting from column J (the last existing data column) into the
# new column K so it matches the rest of the table (fonts, fill, number
# formats, alignment, borders) before adding the right-edge border.
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New 2023 figures for the Gurjaani employment/wages table.
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 3038
$ws.Range("K5").Value = 1137
$ws.Range("K6").Value = 1901

# Column K now closes the table on the right, so give it a thin right
# border matching the rest of the grid.
$ws.Range("K3:K6").Borders.Item(10).LineStyle = 1  # xlEdgeRight, xlContinuous
$ws.Range("K3:K6").Borders.Item(10).Weight = 2     # xlThin

# Match the column width used across the rest of the data columns, and
# widen a few spare columns to the right as well (mirrors the wider
# paintbrush stroke used when the source table was extended further).
$ws.Range("K1:N1").EntireColumn.ColumnWidth = 7.8
